$d = $word.ActiveDocument
$sec = $d.Sections(1)

# There are two logo pictures reused across the document's headers and
# footers (the BTEC logo in both headers, the Pearson logo in both
# footers). Word's own "default name" bookkeeping for these inline
# pictures got out of sync with the actual embedded media part names, so
# rename each inline picture to match what Word now considers correct:
#   headers: BTec logo inline picture -> "image2.jpg"
#   footers: Pearson logo inline picture -> "image1.png"

for ($i = 1; $i -le 2; $i++) {
    $hf = $sec.Headers($i)
    if ($hf.Exists) {
        $shapes = $hf.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shape = $shapes.Item($j)
            if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
                $shape.Name = "image2.jpg"
            }
        }
    }
}

for ($i = 1; $i -le 2; $i++) {
    $ft = $sec.Footers($i)
    if ($ft.Exists) {
        $shapes = $ft.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shape = $shapes.Item($j)
            if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shape.Name = "image1.png"
            }
        }
    }
}
